# Added Test Case and Data for NoOverPay.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PayNowData")

$ws.Range("B7").Value = "6"
$ws.Range("C7").Value = "3.0"
$ws.Range("D7").Value = "10.50"
$ws.Range("E7").Value = "26412171"
$ws.Range("F7").Value = "PayNow"
$ws.Range("G7").Value = "en_US"
$ws.Range("H7").Value = "Elizath"
$ws.Range("I7").Value = "Christine"
$ws.Range("J7").Value = "258 Underwood rd"
$ws.Range("K7").Value = "Suite 600"
$ws.Range("L7").Value = "840"
$ws.Range("M7").Value = "Arlington"
$ws.Range("N7").Value = "VA"
$ws.Range("O7").Value = "22201"

$ws.Range("R7").Value = "iahmed@govolution.com"
$ws.Range("R7").Style = "Normal"
$ws.Range("R7").Borders.LineStyle = 1

$ws.Range("S7").Value = "udf data 1"
$ws.Range("T7").Value = "udf data 2"
$ws.Range("U7").Value = "udf data 3"
$ws.Range("V7").Value = "udf data 4"
$ws.Range("W7").Value = "udf data 5"

$ws.Range("C7").Select()
